# Add a new "About" sheet at the very front of the workbook, describing
# the framework (Name / Description), matching the commit
# "Add about sheet from which name is drawn".

$wb = $excel.ActiveWorkbook

# Insert the new sheet before the current first tab ("Databook Pages")
# so it becomes the new first sheet in the workbook.
$firstSheet = $wb.Worksheets.Item(1)
$aboutSheet = $wb.Worksheets.Add($firstSheet)
$aboutSheet.Name = "About"

# Header row
$aboutSheet.Range("A1").Value = "Name"
$aboutSheet.Range("B1").Value = "Description"
$aboutSheet.Range("A1:B1").Font.Bold = $true

# Content row describing this particular framework
$aboutSheet.Range("A2").Value = "USDT"
$aboutSheet.Range("B2").Value = "USDT cascade"
$aboutSheet.Range("A2:B2").VerticalAlignment = -4160  # xlVAlignTop

# Leave the same selection/active cell Excel would land on after typing
# the content above (matches the saved selection state).
$aboutSheet.Range("B3").Select() | Out-Null
